# Retraining the model for Horeco
# Update Consumption_Actual data: shift dates forward by 2 days and extend range to row 47.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the newly-added rows (34:47) use the same number format style as
# the existing timestamp column (B2:B33, numFmtId 164 "YYYY-MM-DD HH:MM:SS").
$ws.Range("B34:B47").NumberFormat = $ws.Range("B33").NumberFormat

$ws.Cells.Item(2, 1).Value = 5874
$ws.Cells.Item(2, 2).Value = 46047.95833333334
$ws.Cells.Item(3, 1).Value = 5844
$ws.Cells.Item(3, 2).Value = 46047.96875
$ws.Cells.Item(4, 1).Value = 5831
$ws.Cells.Item(4, 2).Value = 46047.97916666666
$ws.Cells.Item(5, 1).Value = 5746
$ws.Cells.Item(5, 2).Value = 46047.98958333334
$ws.Cells.Item(6, 1).Value = 5723
$ws.Cells.Item(6, 2).Value = 46048
$ws.Cells.Item(7, 1).Value = 5712
$ws.Cells.Item(7, 2).Value = 46048.01041666666
$ws.Cells.Item(8, 1).Value = 5692
$ws.Cells.Item(8, 2).Value = 46048.02083333334
$ws.Cells.Item(9, 1).Value = 5629
$ws.Cells.Item(9, 2).Value = 46048.03125
$ws.Cells.Item(10, 1).Value = 5598
$ws.Cells.Item(10, 2).Value = 46048.04166666666
$ws.Cells.Item(11, 1).Value = 5619
$ws.Cells.Item(11, 2).Value = 46048.05208333334
$ws.Cells.Item(12, 1).Value = 5596
$ws.Cells.Item(12, 2).Value = 46048.0625
$ws.Cells.Item(13, 1).Value = 5616
$ws.Cells.Item(13, 2).Value = 46048.07291666666
$ws.Cells.Item(14, 1).Value = 5596
$ws.Cells.Item(14, 2).Value = 46048.08333333334
$ws.Cells.Item(15, 1).Value = 5593
$ws.Cells.Item(15, 2).Value = 46048.09375
$ws.Cells.Item(16, 1).Value = 5590
$ws.Cells.Item(16, 2).Value = 46048.10416666666
$ws.Cells.Item(17, 1).Value = 5618
$ws.Cells.Item(17, 2).Value = 46048.11458333334
$ws.Cells.Item(18, 1).Value = 5651
$ws.Cells.Item(18, 2).Value = 46048.125
$ws.Cells.Item(19, 1).Value = 5650
$ws.Cells.Item(19, 2).Value = 46048.13541666666
$ws.Cells.Item(20, 1).Value = 5686
$ws.Cells.Item(20, 2).Value = 46048.14583333334
$ws.Cells.Item(21, 1).Value = 5725
$ws.Cells.Item(21, 2).Value = 46048.15625
$ws.Cells.Item(22, 1).Value = 5821
$ws.Cells.Item(22, 2).Value = 46048.16666666666
$ws.Cells.Item(23, 1).Value = 5781
$ws.Cells.Item(23, 2).Value = 46048.17708333334
$ws.Cells.Item(24, 1).Value = 5892
$ws.Cells.Item(24, 2).Value = 46048.1875
$ws.Cells.Item(25, 1).Value = 5978
$ws.Cells.Item(25, 2).Value = 46048.19791666666
$ws.Cells.Item(26, 1).Value = 6299
$ws.Cells.Item(26, 2).Value = 46048.20833333334
$ws.Cells.Item(27, 1).Value = 6494
$ws.Cells.Item(27, 2).Value = 46048.21875
$ws.Cells.Item(28, 1).Value = 6678
$ws.Cells.Item(28, 2).Value = 46048.22916666666
$ws.Cells.Item(29, 1).Value = 6973
$ws.Cells.Item(29, 2).Value = 46048.23958333334
$ws.Cells.Item(30, 1).Value = 7330
$ws.Cells.Item(30, 2).Value = 46048.25
$ws.Cells.Item(31, 1).Value = 7621
$ws.Cells.Item(31, 2).Value = 46048.26041666666
$ws.Cells.Item(32, 1).Value = 7735
$ws.Cells.Item(32, 2).Value = 46048.27083333334
$ws.Cells.Item(33, 1).Value = 7886
$ws.Cells.Item(33, 2).Value = 46048.28125
$ws.Cells.Item(34, 1).Value = 8117
$ws.Cells.Item(34, 2).Value = 46048.29166666666
$ws.Cells.Item(35, 1).Value = 8228
$ws.Cells.Item(35, 2).Value = 46048.30208333334
$ws.Cells.Item(36, 1).Value = 8411
$ws.Cells.Item(36, 2).Value = 46048.3125
$ws.Cells.Item(37, 1).Value = 8485
$ws.Cells.Item(37, 2).Value = 46048.32291666666
$ws.Cells.Item(38, 1).Value = 8566
$ws.Cells.Item(38, 2).Value = 46048.33333333334
$ws.Cells.Item(39, 1).Value = 8628
$ws.Cells.Item(39, 2).Value = 46048.34375
$ws.Cells.Item(40, 1).Value = 8673
$ws.Cells.Item(40, 2).Value = 46048.35416666666
$ws.Cells.Item(41, 1).Value = 8711
$ws.Cells.Item(41, 2).Value = 46048.36458333334
$ws.Cells.Item(42, 1).Value = 8627
$ws.Cells.Item(42, 2).Value = 46048.375
$ws.Cells.Item(43, 1).Value = 8617
$ws.Cells.Item(43, 2).Value = 46048.38541666666
$ws.Cells.Item(44, 1).Value = 8558
$ws.Cells.Item(44, 2).Value = 46048.39583333334
$ws.Cells.Item(45, 1).Value = 8571
$ws.Cells.Item(45, 2).Value = 46048.40625
$ws.Cells.Item(46, 1).Value = 8556
$ws.Cells.Item(46, 2).Value = 46048.41666666666
$ws.Cells.Item(47, 1).Value = 8525
$ws.Cells.Item(47, 2).Value = 46048.42708333334
